$d = $word.ActiveDocument

# 1. Add "Checks for invalid geometry" before "Checks field values against
#    template domains where appropriate" (same ListParagraph / numId=4 style).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Checks field values against template domains*") {
        $p.Range.InsertBefore("Checks for invalid geometry`r")
        break
    }
}

# 2. Add "Checks road address range directionality" after "Checks if road
#    features have any geometry cutbacks" (same ListParagraph / numId=4 style).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Checks if road features have any geometry cutbacks*") {
        $p.Range.InsertAfter("`rChecks road address range directionality")
        break
    }
}

# 3. Remove the (hidden) "_GoBack" bookmark that sits right after
#    "ClearOldResults".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 4. Collapse the run of three empty paragraphs (right after "...text files
#    of resource information") down to just one.
$emptyParas = @()
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "`r") {
        $emptyParas += $p
    }
}
$group = @()
foreach ($p in $emptyParas) {
    $prevText = $p.Previous(1).Range.Text
    if ($prevText -like "*text files of resource information*") {
        $group += $p
    }
}
if ($group.Count -gt 0) {
    $first = $group[0]
    $next1 = $first.Next(1)
    $next2 = $next1.Next(1)
    $next2.Range.Delete() | Out-Null
    $next1.Range.Delete() | Out-Null
}

# 5. Split the disclaimer run right after "...In no e" and move the
#    "_GoBack" bookmark into that split point (so it now sits inside the
#    disclaimer text instead of after "ClearOldResults").
$splitRange = $d.Content
$splitRange.Find.Execute("disclaimed.  In no e", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$splitPoint = $d.Range($splitRange.End, $splitRange.End)
$d.Bookmarks.Add("_GoBack", $splitPoint) | Out-Null

Write-Output "done"
